$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The card width/height ("delta" row 5) are now sourced from the actual SVG
# object's measured size instead of hand-rounded numbers, and the card
# border offset (row 24) was re-measured too. Every other D/E cell on the
# sheet is a formula that derives its value from these three inputs, so
# updating them alone ripples through the whole layout table.
$ws.Range("D5").Value = 36.3
$ws.Range("E5").Value = 54.3
$ws.Range("D24").Value = 80

# Reflect the scrolled/selected state of the sheet at save time: the frozen
# bottom pane is scrolled up to row 10 and the active cell moved from N28 to
# D28.
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("A4").Select()
$ws.Range("D28").Select()
